# Applies the "Added most of densenet results" commit to the codecarbon
# results sheet:
#   1. Row 29 was a blank separator row that had been populated with
#      explicit (but empty) string cells -- clear it back to a true blank
#      row while leaving every other row's numbering untouched.
#   2. Four new DenseNet result blocks (CIFAR10 x3, MNIST x2, separated by
#      blank rows, mirroring the existing ResNet blocks above them) are
#      appended as rows 44-112.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Blank out the stray separator row at A29:B29 ---
$ws.Range("A29:B29").ClearContents()

# --- 2. Append the new DenseNet result rows (44-112) ---
$ws.Range("A44").Value = 'Model'
$ws.Range("B44").Value = 'DenseNet'
$ws.Range("A45").Value = 'Dataset'
$ws.Range("B45").Value = 'CIFAR10'
$ws.Range("A46").Value = 'Evaluation Framework'
$ws.Range("B46").Value = 'codecarbon'
$ws.Range("A47").Value = 'Total Energy (kWh)'
$ws.Range("B47").Value = 0.007889708648527477
$ws.Range("A48").Value = 'Total CO2 Emissions (kgCO2e)'
$ws.Range("B48").Value = 0.003665334359850078
$ws.Range("A49").Value = 'CPU Energy'
$ws.Range("B49").Value = 0.003676336044128321
$ws.Range("A50").Value = 'GPU Energy'
$ws.Range("B50").Value = 0.003868527261486
$ws.Range("A51").Value = 'Training Time (minutes)'
$ws.Range("B51").Value = 3.529349748293559
$ws.Range("A52").Value = 'Accuracy'
$ws.Range("B52").Value = 0.4979
$ws.Range("A53").Value = 'Precision'
$ws.Range("B53").Value = 0.5076567591666675
$ws.Range("A54").Value = 'Recall'
$ws.Range("B54").Value = 0.4979
$ws.Range("A55").Value = 'F1'
$ws.Range("B55").Value = 0.4938883820971905
$ws.Range("A56").Value = 'Number of Epochs'
$ws.Range("B56").Value = 10
$ws.Range("A58").Value = 'Model'
$ws.Range("B58").Value = 'DenseNet'
$ws.Range("A59").Value = 'Dataset'
$ws.Range("B59").Value = 'CIFAR10'
$ws.Range("A60").Value = 'Evaluation Framework'
$ws.Range("B60").Value = 'codecarbon'
$ws.Range("A61").Value = 'Total Energy (kWh)'
$ws.Range("B61").Value = 0.01376071803926909
$ws.Range("A62").Value = 'Total CO2 Emissions (kgCO2e)'
$ws.Range("B62").Value = 0.006392838429459039
$ws.Range("A63").Value = 'CPU Energy'
$ws.Range("B63").Value = 0.006502480990958347
$ws.Range("A64").Value = 'GPU Energy'
$ws.Range("B64").Value = 0.006648244763035999
$ws.Range("A65").Value = 'Training Time (minutes)'
$ws.Range("B65").Value = 6.24258105357488
$ws.Range("A66").Value = 'Accuracy'
$ws.Range("B66").Value = 0.6074000000000001
$ws.Range("A67").Value = 'Precision'
$ws.Range("B67").Value = 0.60533981406285
$ws.Range("A68").Value = 'Recall'
$ws.Range("B68").Value = 0.6074000000000001
$ws.Range("A69").Value = 'F1'
$ws.Range("B69").Value = 0.6015483504098486
$ws.Range("A70").Value = 'Number of Epochs'
$ws.Range("B70").Value = 10
$ws.Range("A72").Value = 'Model'
$ws.Range("B72").Value = 'DenseNet'
$ws.Range("A73").Value = 'Dataset'
$ws.Range("B73").Value = 'CIFAR10'
$ws.Range("A74").Value = 'Evaluation Framework'
$ws.Range("B74").Value = 'codecarbon'
$ws.Range("A75").Value = 'Total Energy (kWh)'
$ws.Range("B75").Value = 0.01381779868302146
$ws.Range("A76").Value = 'Total CO2 Emissions (kgCO2e)'
$ws.Range("B76").Value = 0.00641935647393296
$ws.Range("A77").Value = 'CPU Energy'
$ws.Range("B77").Value = 0.006506989855818876
$ws.Range("A78").Value = 'GPU Energy'
$ws.Range("B78").Value = 0.006700420915888
$ws.Range("A79").Value = 'Training Time (minutes)'
$ws.Range("B79").Value = 6.246892023086548
$ws.Range("A80").Value = 'Accuracy'
$ws.Range("B80").Value = 0.8308
$ws.Range("A81").Value = 'Precision'
$ws.Range("B81").Value = 0.8312846509078877
$ws.Range("A82").Value = 'Recall'
$ws.Range("B82").Value = 0.8308
$ws.Range("A83").Value = 'F1'
$ws.Range("B83").Value = 0.8308816365228535
$ws.Range("A84").Value = 'Number of Epochs'
$ws.Range("B84").Value = 10
$ws.Range("A86").Value = 'Model'
$ws.Range("B86").Value = 'DenseNet'
$ws.Range("A87").Value = 'Dataset'
$ws.Range("B87").Value = 'MNIST'
$ws.Range("A88").Value = 'Evaluation Framework'
$ws.Range("B88").Value = 'codecarbon'
$ws.Range("A89").Value = 'Total Energy (kWh)'
$ws.Range("B89").Value = 0.02369511978299916
$ws.Range("A90").Value = 'Total CO2 Emissions (kgCO2e)'
$ws.Range("B90").Value = 0.01100807907749544
$ws.Range("A91").Value = 'CPU Energy'
$ws.Range("B91").Value = 0.01102109265664694
$ws.Range("A92").Value = 'GPU Energy'
$ws.Range("B92").Value = 0.01164026486776
$ws.Range("A93").Value = 'Training Time (minutes)'
$ws.Range("B93").Value = 10.58042073647181
$ws.Range("A94").Value = 'Accuracy'
$ws.Range("B94").Value = 0.991
$ws.Range("A95").Value = 'Precision'
$ws.Range("B95").Value = 0.9910263475630446
$ws.Range("A96").Value = 'Recall'
$ws.Range("B96").Value = 0.991
$ws.Range("A97").Value = 'F1'
$ws.Range("B97").Value = 0.9909997736228469
$ws.Range("A98").Value = 'Number of Epochs'
$ws.Range("B98").Value = 1

# Row 99 is a blank separator row that still carries explicit empty-string
# cells (matches the pattern used elsewhere in this sheet, e.g. old row 29).
# Assign an empty-string formula so the engine commits real (empty) text
# cells here instead of leaving the row completely absent.
$ws.Range("A99").Formula = '=""'
$ws.Range("B99").Formula = '=""'

$ws.Range("A100").Value = 'Model'
$ws.Range("B100").Value = 'DenseNet'
$ws.Range("A101").Value = 'Dataset'
$ws.Range("B101").Value = 'MNIST'
$ws.Range("A102").Value = 'Evaluation Framework'
$ws.Range("B102").Value = 'codecarbon'
$ws.Range("A103").Value = 'Total Energy (kWh)'
$ws.Range("B103").Value = 0.2368809560385094
$ws.Range("A104").Value = 'Total CO2 Emissions (kgCO2e)'
$ws.Range("B104").Value = 0.110048158435373
$ws.Range("A105").Value = 'CPU Energy'
$ws.Range("B105").Value = 0.1101841889024743
$ws.Range("A106").Value = 'GPU Energy'
$ws.Range("B106").Value = 0.116361546978052
$ws.Range("A107").Value = 'Training Time (minutes)'
$ws.Range("B107").Value = 105.777428428332
$ws.Range("A108").Value = 'Accuracy'
$ws.Range("B108").Value = 0.9957
$ws.Range("A109").Value = 'Precision'
$ws.Range("B109").Value = 0.9957108101984441
$ws.Range("A110").Value = 'Recall'
$ws.Range("B110").Value = 0.9957
$ws.Range("A111").Value = 'F1'
$ws.Range("B111").Value = 0.9956985806654162
$ws.Range("A112").Value = 'Number of Epochs'
$ws.Range("B112").Value = 10
